# Add 2022-Q3 data
# 1) Insert a new worksheet "2022-Q3" right after "总计" and before "2022-Q2",
#    by copying the existing "2022-Q2" sheet (to inherit all formatting/styles)
#    and then overwriting its data with the Q3 figures.
# 2) Insert a new row into the "总计" (summary) sheet for the 2022-Q3 totals,
#    shifting the existing 2022-Q2 / 2022-Q1 / 2021-Q4 rows down by one.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Step 1: build the new "2022-Q3" sheet
# ---------------------------------------------------------------------------
$templateSheet = $wb.Worksheets.Item("2022-Q2")
$templateSheet.Copy($templateSheet) | Out-Null

$q3 = $wb.Worksheets.Item(2)
$q3.Name = "2022-Q3"

# The template ("2022-Q2") had 8 data rows (rows 2-9); the new quarter only
# has 2 funds, so remove the now-unused rows 4-9.
$q3.Range("A4:H9").Delete() | Out-Null

# Force the code / numeric-looking text columns to remain text so that
# leading zeros (fund codes) and decimal formatting are preserved exactly.
$q3.Range("B2:B3").NumberFormat = "@"
$q3.Range("D2:G3").NumberFormat = "@"

$q3.Cells.Item(2, 1).Value = 0
$q3.Cells.Item(2, 2).Value = "010874"
$q3.Cells.Item(2, 3).Value = "泰康品质生活混合A"
$q3.Cells.Item(2, 4).Value = "6.86"
$q3.Cells.Item(2, 5).Value = "84.42"
$q3.Cells.Item(2, 6).Value = "3.26"
$q3.Cells.Item(2, 7).Value = "0.2236"
$q3.Cells.Item(2, 8).Value = 9

$q3.Cells.Item(3, 1).Value = 1
$q3.Cells.Item(3, 2).Value = "010875"
$q3.Cells.Item(3, 3).Value = "泰康品质生活混合C"
$q3.Cells.Item(3, 4).Value = "3.44"
$q3.Cells.Item(3, 5).Value = "84.42"
$q3.Cells.Item(3, 6).Value = "3.26"
$q3.Cells.Item(3, 7).Value = "0.1121"
$q3.Cells.Item(3, 8).Value = 9

# ---------------------------------------------------------------------------
# Step 2: update the "总计" (summary) sheet
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

$total.Rows.Item(2).Insert() | Out-Null

# The "A" column is a running 0-based index; bump the (now shifted down)
# existing rows by one and give the new row index 0.
$total.Cells.Item(5, 1).Value = 3
$total.Cells.Item(4, 1).Value = 2
$total.Cells.Item(3, 1).Value = 1

$a2 = $total.Cells.Item(2, 1)
$a2.Value = 0
$a2.Font.Bold = $true
$a2.HorizontalAlignment = -4108
$a2.VerticalAlignment = -4160
$a2.Borders.Item(7).LineStyle = 1
$a2.Borders.Item(8).LineStyle = 1
$a2.Borders.Item(9).LineStyle = 1
$a2.Borders.Item(10).LineStyle = 1

$total.Range("B2:D2").ClearFormats() | Out-Null
$total.Cells.Item(2, 2).Value = "2022-Q3"
$total.Cells.Item(2, 3).Value = 2
$total.Cells.Item(2, 4).Value = 0.34

# Keep the "总计" sheet as the active tab (same as the original workbook).
$total.Activate()
$total.Range("A1").Select() | Out-Null
